$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.314.40"
$ws.Range("E2").Value = "  +3.60%  "

$ws.Range("D3").Value = "2.069.16"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.42%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.89%  "

$ws.Range("E10").Value = "  +5.00%  "

$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.09%  "

$ws.Range("D13").Value = "2.375.15"
$ws.Range("E13").Value = "  +3.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.92%  "

$ws.Range("D17").Value = "2.069.65"
$ws.Range("E17").Value = "  +1.95%  "

$ws.Range("D18").Value = "38.242.54"

$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.68%  "

$ws.Range("E21").Value = "  +3.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "

$ws.Range("E25").Value = "  +4.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.07%  "

$ws.Range("E31").Value = "  +3.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.54%  "

$ws.Range("E33").Value = "  +5.57%  "

$ws.Range("E34").Value = "  +10.00%  "

$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.47%  "

$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("D40").Value = "1.530.11"
$ws.Range("E40").Value = "  +4.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.70%  "

$ws.Range("E43").Value = "  +3.97%  "

$ws.Range("E44").Value = "  +4.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0933"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.62%  "

$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.43%  "

$ws.Range("E48").Value = "  +3.43%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.13%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").Value = "2.262.00"
$ws.Range("E51").Value = "  +3.17%  "
